# Generate Report for Handoff
#
# The "b.md" row (row 3) across the Overview / zh-cn / de-de sheets is
# updated to reflect that the file is ready for handoff: the newly
# generated target xliff ("b.*.xlf") is recorded, along with its
# generation timestamp, and an error noting the handback file used for
# translation is stale compared to the newest source.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet: row 3 is the "b.md" file.
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-27 04:36:27"

# ---------------------------------------------------------------------
# zh-cn sheet: row 3 is the "b.md" file.
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"

# "False" needs to land as literal text (matching the rest of the sheet,
# which stores True/False as plain strings) rather than a native Excel
# boolean, so round-trip it through a text formula + paste-as-values.
$wsZhCn.Range("F3").Formula = "=""False"""
$wsZhCn.Range("F3").Copy()
$wsZhCn.Range("F3").PasteSpecial(-4163)

$wsZhCn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-08-27 04:36:23"
$wsZhCn.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/28a3f882d89c73d7083775a34a5c25f110aefdb6/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9025cfc53cd959d33e1f618583080d48230aae7b/e2e/b.md."
$wsZhCn.Columns.Item(16).ColumnWidth = 39.166666666666664

# ---------------------------------------------------------------------
# de-de sheet: row 3 is the "b.md" file.
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"

$wsDeDe.Range("F3").Formula = "=""False"""
$wsDeDe.Range("F3").Copy()
$wsDeDe.Range("F3").PasteSpecial(-4163)

$wsDeDe.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-08-27 04:36:27"
$wsDeDe.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/28a3f882d89c73d7083775a34a5c25f110aefdb6/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9025cfc53cd959d33e1f618583080d48230aae7b/e2e/b.md."
$wsDeDe.Columns.Item(16).ColumnWidth = 39.166666666666664
